$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---------------
# This string is shared by the Overview sheet's per-language status columns
# (E2 = zh-cn status, F2 = de-de status) and each language sheet's own
# "Status" column (C2). Updating all of them lands on the same shared string.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the "Status" columns --------------------------------------
# Overview: columns E (zh-cn) and F (de-de) both hold status values.
# ColumnWidth = 12.5 is the COM input that this host's pixel-grid rounding
# resolves to the narrower stored width used for the status columns.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C is the "Status" column.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
